$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1312794072705492
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2696233515299001
$ws.Range("E2").Value = 0.007675740627754522
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.09162142013456404
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2092013713451965
$ws.Range("N2").Value = 0.02010645451059558
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1042507738322135
$ws.Range("V2").Value = 0.02676475543772707
$ws.Range("W2").Value = -0.006624392836992566
$ws.Range("Z2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.01229826710515999
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.01603021089254992
$ws.Range("AF2").Value = -0.002920549678736569
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.03213058706619464
$ws.Range("AN2").Value = 0.04866463345903858
$ws.Range("AO2").Value = 0.05872192759438443
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AS2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1420381477075912
$ws.Range("AW2").Value = 0.09604362001676886
$ws.Range("AX2").Value = -0.03007136911367418
$ws.Range("AY2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.01394801060804107
$ws.Range("BF2").Value = 0.09774606981162738
$ws.Range("BG2").Value = 0.01691949969733672
$ws.Range("BI2").Value = 0
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.01336856832920403
$ws.Range("BO2").Value = -0.03872286851754293
$ws.Range("BP2").Value = -0.07789378838805426
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.06541867347822705
$ws.Range("BX2").Value = 0.0286815084985256
$ws.Range("BY2").Value = -0.03667635841468975
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03131965819972365
$ws.Range("CF2").Value = -0
$ws.Range("CG2").Value = -0.03792856407418337
$ws.Range("CH2").Value = 0.01977921938127206
$ws.Range("CI2").Value = 0
$ws.Range("CJ2").Value = -0
$ws.Range("CK2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.05526154718428198
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.04082519066560224
$ws.Range("CQ2").Value = 0.07604816541436506
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04869620284243264
$ws.Range("CY2").Value = -0.04133256821757724
$ws.Range("CZ2").Value = 0.01431381162547751
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.04211323188845133
$ws.Range("DH2").Value = -0.005094718869252873
$ws.Range("DI2").Value = 0.009147846176215805
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01458424503963165
$ws.Range("DQ2").Value = 0.06318707184845157
$ws.Range("DR2").Value = -0.01288648389114542
$ws.Range("DS2").Value = -0
$ws.Range("DT2").Value = 0
$ws.Range("DU2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05106129249218446
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = 0.002891590206178732
$ws.Range("EA2").Value = -0.02682700005420498
$ws.Range("EB2").Value = 0
$ws.Range("ED2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.03505717613743062
$ws.Range("EI2").Value = 0.09774273743666885
$ws.Range("EJ2").Value = -0.01855732252884039
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.06391125286757442
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.05342716592926169
$ws.Range("ES2").Value = 0.01560759277549106
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.03939709954443905
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.0334448179641134
$ws.Range("FB2").Value = 0.02489131420894644
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.002947118385050842
$ws.Range("FJ2").Value = -0.03186331621889533
$ws.Range("FK2").Value = 0.01550483445938226
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = 0.001061672145293953
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = 0.01328256027017473
$ws.Range("FT2").Value = 0.002496242915856067
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03271345762496612
$ws.Range("GB2").Value = 0.02005732247023465
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
